# Apply the commit's changes to the "contact" worksheet.
#
# Summary of the change:
#  - Row 2 (the "tett" / empty / 2022-03-17T09:14:59.586Z / ... record) is
#    updated in place to "tt" / "11" / 2022-03-18T20:02:40.062Z /
#    6234e56efa1a37ad6b10fbb8 / 6234e57afa1a37ad6b10fbc1.
#  - Row 3 is replaced with a brand new record: fawzi / ttt / ttt / 1 /
#    2022-03-18T20:27:31.009Z / 6234e56efa1a37ad6b10fbb8 /
#    6234eb4ed2906ab295e2fb29.
#  - The old row 4 is removed entirely (the sheet now only spans A1:G3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Mobile" column values ("11" and "1") look numeric, so Excel would
# normally auto-convert them to numbers on assignment. Force them to be
# stored as text (matching the original shared-string layout) by
# temporarily using a text number format, then clear the format again so
# the cell keeps using the sheet's default style.
$ws.Range("D2:D3").NumberFormat = "@"

# --- Update row 2 in place ---
$ws.Range("A2").Value = "tt"
$ws.Range("B2").Value = "tt"
$ws.Range("C2").Value = "tt"
$ws.Range("D2").Value = "11"
$ws.Range("E2").Value = "2022-03-18T20:02:40.062Z"
$ws.Range("F2").Value = "6234e56efa1a37ad6b10fbb8"
$ws.Range("G2").Value = "6234e57afa1a37ad6b10fbc1"

# --- Replace row 3 with the new record ---
$ws.Range("A3").Value = "fawzi"
$ws.Range("B3").Value = "ttt"
$ws.Range("C3").Value = "ttt"
$ws.Range("D3").Value = "1"
$ws.Range("E3").Value = "2022-03-18T20:27:31.009Z"
$ws.Range("F3").Value = "6234e56efa1a37ad6b10fbb8"
$ws.Range("G3").Value = "6234eb4ed2906ab295e2fb29"

# Drop the temporary text format again so the cells fall back to the
# workbook's default (unstyled) cell style, exactly like every other cell.
$ws.Range("D2:D3").ClearFormats()

# --- Remove the old row 4 entirely ---
$ws.Rows.Item(4).Delete()
